$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A..R) updates
$ws.Range("A2").Value = 112253960
$ws.Range("B2").Value = 90466
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = "Svavelriska"
$ws.Range("G2").Value = "Lactarius scrobiculatus"
$ws.Range("H2").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q2").Value = 571434
$ws.Range("R2").Value = 6703294

# Row 3 (A..R) updates
$ws.Range("A3").Value = 112253958
$ws.Range("B3").Value = 98961
$ws.Range("E3").Value = 222498
$ws.Range("F3").Value = "Blåsippa"
$ws.Range("G3").Value = "Hepatica nobilis"
$ws.Range("H3").Value = "Schreb."
$ws.Range("Q3").Value = 571485
$ws.Range("R3").Value = 6703317
